$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row 1: add P1=14 and Q1=15, matching the style of O1 ---
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- Update data rows 2-25 ---
for ($r = 2; $r -le 25; $r++) {
    # Swap values: I<->K style swap (1<->2), M<->O style swap (1<->2)
    $ws.Cells.Item($r, 9).Value = 2   # I: 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1  # K: 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2  # M: 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1  # O: 2 -> 1

    # Add new columns P and Q with value 2 (no special style, like other data cells)
    $ws.Cells.Item($r, 16).Value = 2  # P
    $ws.Cells.Item($r, 17).Value = 2  # Q
}
